$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.784.46'
$ws.Range("E2").Value = '  -0.33%  '
$ws.Range("D3").Value = '3.307.67'
$ws.Range("E3").Value = '  +0.73%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.64'
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '185.29'
$ws.Range("E6").Value = '  -0.29%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.298.98'
$ws.Range("E8").Value = '  +0.64%  '
$ws.Range("E9").Value = '  -3.33%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.175'
$ws.Range("E10").Value = '  -6.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.576'
$ws.Range("E11").Value = '  -1.97%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.63'
$ws.Range("E12").Value = '  -4.05%  '
$ws.Range("E13").Value = '  -2.37%  '
$ws.Range("D14").Value = '3.838.72'
$ws.Range("E14").Value = '  +0.63%  '
$ws.Range("E15").Value = '  -2.70%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '572.90'
$ws.Range("E16").Value = '  -9.94%  '
$ws.Range("D17").Value = '65.722.23'
$ws.Range("E17").Value = '  -0.55%  '
$ws.Range("E18").Value = '  +0.33%  '
$ws.Range("D19").Value = '3.304.12'
$ws.Range("E19").Value = '  +0.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.62'
$ws.Range("E20").Value = '  -2.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.77'
$ws.Range("E21").Value = '  -5.14%  '
$ws.Range("E22").Value = '  -2.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.86'
$ws.Range("E23").Value = '  -2.85%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.98'
$ws.Range("E24").Value = '  +1.62%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '98.49'
$ws.Range("E25").Value = '  -8.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.92'
$ws.Range("E26").Value = '  -0.98%  '
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.29'
$ws.Range("E28").Value = '  -3.18%  '
$ws.Range("E29").Value = '  -3.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.41'
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.62'
$ws.Range("E31").Value = '  +5.75%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.67'
$ws.Range("E32").Value = '  -8.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '556.16'
$ws.Range("E33").Value = '  +5.44%  '
$ws.Range("E34").Value = '  -2.34%  '
$ws.Range("D35").Value = '3.760.91'
$ws.Range("E35").Value = '  +1.45%  '
$ws.Range("E36").Value = '  -2.72%  '
$ws.Range("E37").Value = '  +0.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.59'
$ws.Range("E38").Value = '  -3.62%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '33.55'
$ws.Range("E39").Value = '  +1.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.126'
$ws.Range("E40").Value = '  -3.24%  '
$ws.Range("D41").Value = '0.0₃0678'
$ws.Range("E41").Value = '  -6.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.10'
$ws.Range("E42").Value = '  -7.37%  '
$ws.Range("E45").Value = '  -2.65%  '
$ws.Range("E46").Value = '  -2.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.00'
$ws.Range("E47").Value = '  -11.37%  '
$ws.Range("E50").Value = '  -3.78%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '124.90'
$ws.Range("E51").Value = '  +2.67%  '

$ws.Range("B43").Value = 'ApeXProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.33'
$ws.Range("E43").Value = '  +1.31%  '
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.56'
$ws.Range("E44").Value = '  -5.83%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.126'
$ws.Range("E48").Value = '  -2.45%  '
$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.999'
$ws.Range("E49").Value = '  -0.06%  '
